$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list: refresh Price (column D) and Volume(1h) (column E)
# values for this data pull. Values are stored as text in the sheet, so we
# force a text number-format while writing, then restore the cell style to
# "Normal" so no stray formatting is left behind.

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Value
    )
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "243.15"
Set-TextValue "D3" "23.06"
Set-TextValue "D4" "5.400"
Set-TextValue "D5" "0.05921"
Set-TextValue "D6" "3.449"
Set-TextValue "D7" "6.529"
Set-TextValue "D8" "0.8109"
Set-TextValue "D10" "0.1407"
Set-TextValue "D11" "0.07376"
Set-TextValue "D12" "0.03269"
Set-TextValue "D13" "0.03044"
Set-TextValue "D14" "0.09355"
Set-TextValue "D15" "3.849"
Set-TextValue "D16" "0.001583"
Set-TextValue "D17" "0.04677"
Set-TextValue "D18" "0.0005940"
Set-TextValue "D19" "0.006094"
Set-TextValue "D20" "0.004977"
Set-TextValue "D21" "0.0009815"
Set-TextValue "D22" "0.00009404"
Set-TextValue "E22" "21NitroExNTXBestin24h"
Set-TextValue "D23" "3.608"
Set-TextValue "D24" "2.148"
Set-TextValue "D40" "0.03964"
Set-TextValue "E41" "40KickTokenKICK"
Set-TextValue "D42" "0.1077"
Set-TextValue "D43" "0.002621"
Set-TextValue "D44" "0.008035"
Set-TextValue "D45" "0.00005255"
Set-TextValue "D48" "0.002265"
